# Reorganize the functions / prompts section at the end of the document:
#   - remove the "My Notes" heading and all of the "// ..." planning-comment
#     paragraphs that followed it
#   - remove the two blank paragraphs that preceded the page break
#   - remove the page-break run itself, while keeping the (now empty) paragraph
#     that used to host it (its paragraph mark / pPr stay in the document)

$d = $word.ActiveDocument

# --- locate the "My Notes" paragraph (first paragraph of the trailing notes block) ---
$rng = $d.Content
$found = $rng.Find.Execute("My Notes", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'My Notes' paragraph"
}
$notesStart = $rng.Start

# Figure out which paragraph index that start offset belongs to.
$notesIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $notesStart -and $notesStart -lt $p.Range.End) {
        $notesIndex = $i
        break
    }
}
if ($notesIndex -eq 0) {
    throw "Could not resolve paragraph index for 'My Notes'"
}

$lastIndex = $d.Paragraphs.Count

# --- delete "My Notes" through the final "// if no, re-run all the generators" paragraph ---
for ($i = $lastIndex; $i -ge $notesIndex; $i--) {
    $d.Paragraphs.Item($i).Range.Delete()
}

# --- the paragraph that now sits right before where "My Notes" used to be is the
#     page-break paragraph; the two paragraphs before that are blank and get removed ---
$pageBreakIndex = $notesIndex - 1
$blankBeforeBreakIndex = $pageBreakIndex - 1
$blankFirstIndex = $pageBreakIndex - 2

$d.Paragraphs.Item($blankBeforeBreakIndex).Range.Delete()
$d.Paragraphs.Item($blankFirstIndex).Range.Delete()

# Both blank paragraphs are now gone, so the page-break paragraph has shifted down
# to where the first blank paragraph used to be.
$pageBreakIndex = $blankFirstIndex

# --- remove the page-break run from that paragraph, leaving the empty paragraph
#     (and its pPr/paragraph mark) intact ---
$breakPara = $d.Paragraphs.Item($pageBreakIndex)
$pStart = $breakPara.Range.Start
$pEnd = $breakPara.Range.End
$runRange = $d.Range($pStart, $pEnd - 1)
if ($runRange.Text.Length -gt 0) {
    $runRange.Delete()
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
